$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "4 Owned" column header to "# Owned"
$ws.Range("C1").Value = "# Owned"

# Reflect the active selection left by the edit (matches the saved file)
$ws.Range("C1").Select()
